$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
}

# Row 2
$ws.Range("B2").Value = "Mapping Global Cyberterror Networks: An Empirical Study of Al-Qaeda and ISIS Cyberterrorism Events"
$ws.Range("C2").Value = "Claire Seungeun Lee, Kyung-Shick Choi, Ryan Shandler, Chris Kayser"
Set-TextValue $ws.Range("D2") "2021"
$ws.Range("E2").Value = "10.1177/10439862211001606"
$ws.Range("G2").Value = 15

# Row 3
$ws.Range("B3").Value = "Responding to Uncertainty: The Importance of Covertness in Support for Retaliation to Cyber and Kinetic Attacks"
$ws.Range("C3").Value = "Kathryn Hedgecock, Lauren Sukin"
Set-TextValue $ws.Range("D3") "2023"
$ws.Range("E3").Value = "10.1177/00220027231153580"
$ws.Range("F3").Value = "Open Access"
$ws.Range("G3").Value = 11

# Row 4
$ws.Range("B4").Value = "Fighting in Cyberspace: Internet Access and the Substitutability of Cyber and Military Operations"
$ws.Range("C4").Value = "Nadiya Kostyuk, Erik Gartzke"
Set-TextValue $ws.Range("D4") "2024"
$ws.Range("E4").Value = "10.1177/00220027231160993"
$ws.Range("G4").Value = 9

# Row 5
$ws.Range("B5").Value = "Attrition rates and maneuver in agent-based simulation models"
$ws.Range("C5").Value = "David Ormrod, Benjamin Turnbull"
Set-TextValue $ws.Range("D5") "2017"
$ws.Range("E5").Value = "10.1177/1548512917692693"
$ws.Range("G5").Value = 17

# Row 6
$ws.Range("B6").Value = "Indonesia’s Handling of Terrorists’ Cyber Activities: How Repressive Measures Still Fall Short"
$ws.Range("C6").Value = "Ali Abdullah Wibisono, Rachel Kumendong, Iwa Maulana"
Set-TextValue $ws.Range("D6") "2025"
$ws.Range("E6").Value = "10.1177/23477970241298764"
$ws.Range("G6").Value = 9

# Row 7
$ws.Range("B7").Value = "Wargaming the use of intermediate force capabilities in the gray zone"
$ws.Range("C7").Value = "Kyle D Christensen, Peter Dobias"
$ws.Range("E7").Value = "10.1177/15485129211010227"
$ws.Range("F7").Value = "Restricted"
$ws.Range("G7").Value = 11

# Row 8
$ws.Range("B8").Value = "Robust tracking strategy for nonlinear connected vehicle cyber-physical systems"
$ws.Range("C8").Value = "Yushi Yang, Meng Li, Yong Chen"
Set-TextValue $ws.Range("D8") "2024"
$ws.Range("E8").Value = "10.1177/01423312231196642"
$ws.Range("F8").Value = "Restricted"
$ws.Range("G8").Value = 19

# Row 9
$ws.Range("B9").Value = "Tech titans, cyber commons and the war in Ukraine: An incipient shift in international relations"
$ws.Range("C9").Value = "Eviatar Matania, Udi Sommer"
Set-TextValue $ws.Range("D9") "2023"
$ws.Range("E9").Value = "10.1177/00471178231211500"
$ws.Range("F9").Value = "Open Access"
$ws.Range("G9").Value = 22

# Row 10
$ws.Range("B10").Value = "Using network digital twins to improve cyber resilience of missions"
$ws.Range("C10").Value = "Rajive Bagrodia"
Set-TextValue $ws.Range("D10") "2023"
$ws.Range("E10").Value = "10.1177/15485129221131226"

# Row 11
$ws.Range("B11").Value = "Simplification and Linearization of Manipulator Dynamics by the Design of Inertia Distribution"
$ws.Range("C11").Value = "D.C.H. Yang, S.W. Tzeng"
Set-TextValue $ws.Range("D11") "1986"
$ws.Range("E11").Value = "10.1177/027836498600500307"
$ws.Range("G11").Value = 4

$wb.Save()
